# Updates cryptos list (prices / 1h volume changes) and re-ranks a few coins
# whose relative order changed (rows 26-29), matching the GitHub Actions
# refresh commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.091.83'
$ws.Range("E2").Value = '  +5.31%  '

$ws.Range("D3").Value = '3.329.39'
$ws.Range("E3").Value = '  +1.88%  '

$ws.Range("E4").Value = '  +0.13%  '

$ws.Range("D5").Value = "'411.31"
$ws.Range("E5").Value = '  +3.56%  '

$ws.Range("D6").Value = "'110.95"
$ws.Range("E6").Value = '  +1.06%  '

$ws.Range("D7").Value = "'0.583"
$ws.Range("E7").Value = '  +3.86%  '

$ws.Range("E8").Value = '  +0.07%  '

$ws.Range("D9").Value = "'0.631"
$ws.Range("E9").Value = '  +0.58%  '

$ws.Range("D10").Value = "'39.32"
$ws.Range("E10").Value = '  +0.01%  '

$ws.Range("D11").Value = "'0.0980"
$ws.Range("E11").Value = '  +1.82%  '

$ws.Range("E12").Value = '  +0.95%  '

$ws.Range("D13").Value = '3.866.54'
$ws.Range("E13").Value = '  +2.40%  '

$ws.Range("D14").Value = "'8.38"
$ws.Range("E14").Value = '  +1.97%  '

$ws.Range("D15").Value = "'19.56"
$ws.Range("E15").Value = '  +2.41%  '

$ws.Range("D16").Value = '3.330.60'
$ws.Range("E16").Value = '  +2.15%  '

$ws.Range("E17").Value = '  +0.01%  '

$ws.Range("D18").Value = '59.919.48'
$ws.Range("E18").Value = '  +5.43%  '

$ws.Range("D19").Value = "'10.72"
$ws.Range("E19").Value = '  -0.91%  '

$ws.Range("D20").Value = "'3.36"
$ws.Range("E20").Value = '  +1.81%  '

$ws.Range("E21").Value = '  +3.49%  '

$ws.Range("E22").Value = '  +1.81%  '

$ws.Range("D23").Value = "'299.04"
$ws.Range("E23").Value = '  -2.79%  '

$ws.Range("D24").Value = "'75.00"
$ws.Range("E24").Value = '  -0.35%  '

$ws.Range("D25").Value = "'3.18"
$ws.Range("E25").Value = '  +0.71%  '

$ws.Range("B26").Value = 'Filecoin'
$ws.Range("C26").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D26").Value = "'8.08"
$ws.Range("E26").Value = '  +2.37%  '

$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").Value = "'28.51"
$ws.Range("E27").Value = '  +1.30%  '

$ws.Range("B28").Value = 'RenderToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D28").Value = "'7.75"
$ws.Range("E28").Value = '  +6.78%  '

$ws.Range("B29").Value = 'LEO'
$ws.Range("C29").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D29").Value = "'4.45"
$ws.Range("E29").Value = '  +1.60%  '

$ws.Range("D30").Value = "'0.180"
$ws.Range("E30").Value = '  +5.82%  '

$ws.Range("D31").Value = "'0.115"
$ws.Range("E31").Value = '  +4.09%  '

$ws.Range("D32").Value = "'2.59"
$ws.Range("E32").Value = '  +20.84%  '

$ws.Range("E33").Value = '  +3.71%  '

$ws.Range("E34").Value = '  +0.24%  '

$ws.Range("D35").Value = "'39.33"
$ws.Range("E35").Value = '  +4.70%  '

$ws.Range("D36").Value = "'0.0503"
$ws.Range("E36").Value = '  +4.30%  '

$ws.Range("D37").Value = "'52.12"
$ws.Range("E37").Value = '  +1.09%  '

$ws.Range("E38").Value = '  +0.18%  '

$ws.Range("D39").Value = "'3.09"
$ws.Range("E39").Value = '  -0.52%  '

$ws.Range("D40").Value = "'3.39"
$ws.Range("E40").Value = '  -3.77%  '

$ws.Range("D41").Value = "'138.02"
$ws.Range("E41").Value = '  +2.19%  '

$ws.Range("E42").Value = '  +2.07%  '

$ws.Range("D43").Value = "'0.293"
$ws.Range("E43").Value = '  +4.11%  '

$ws.Range("E44").Value = '  -0.99%  '

$ws.Range("D45").Value = "'3.91"
$ws.Range("E45").Value = '  -2.35%  '

$ws.Range("D46").Value = "'16.82"
$ws.Range("E46").Value = '  -2.29%  '

$ws.Range("D47").Value = "'2.26"
$ws.Range("E47").Value = '  +8.63%  '

$ws.Range("D48").Value = "'22.15"
$ws.Range("E48").Value = '  -0.07%  '

$ws.Range("D49").Value = '2.185.31'
$ws.Range("E49").Value = '  +1.61%  '

$ws.Range("E50").Value = '  +1.45%  '

$ws.Range("D51").Value = "'2.02"
$ws.Range("E51").Value = '  +0.25%  '
